# Updating calculations, sensors page and adding micro page
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Enable Circuit (VBATT UVLO) section ---
# R1
$ws.Range("B5").Value = 130000
# R2
$ws.Range("B6").Value = 12000
# V(en) (B7) recalculates automatically from =B4*(B6/(B6+B5))

# --- LEDs section ---
# Vf (B27) was a formula (=4.2*4); replace with the plain computed value
$ws.Range("B27").Value = 14.4
# R (B31) recalculates automatically from =B27*(B30/(B30+B29))

# C36 recalculates automatically from =(B27-C34)/C35

# Move the selection to B28, with no frozen/scrolled topLeftCell
$ws.Range("B28").Select()
